$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 17, shifting the existing rows 17-20 down to 18-21.
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with the latest weekly price record.
$ws.Cells.Item(17, 1).Value = 7
$ws.Cells.Item(17, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(17, 3).Value = "Ñuble"
$ws.Cells.Item(17, 4).Value = 44841
$ws.Cells.Item(17, 5).Value = 16
$ws.Cells.Item(17, 6).Value = 100112012
$ws.Cells.Item(17, 7).Value = "Espinaca"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 60
$ws.Cells.Item(17, 11).Value = 6500
$ws.Cells.Item(17, 12).Value = 7000
$ws.Cells.Item(17, 13).Value = 6750
$ws.Cells.Item(17, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(17, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(17, 16).Value = 675
$ws.Cells.Item(17, 17).Value = 10
$ws.Cells.Item(17, 18).Value = "Hortaliza"
